$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-23 (tech / modelTech pairs) are being reordered: "IndustryH_Biomass"
# moves from row 22 up to row 7, and "GT" moves from row 23 up to row 18,
# with the intervening rows shifting down by one position accordingly.
# Rewrite columns A and B for rows 7-23 with the new, reordered values.

$newTech = @(
    "IndustryH_Biomass",
    "BP_Coal",
    "BH_Natgas",
    "BP_Natgas",
    "IndustryH_Natgas",
    "BH_Oil",
    "BP_Oil",
    "IndustryH_Oil",
    "BH_Waste",
    "BP_Waste",
    "EP",
    "GT",
    "HPstandard",
    "HPsurplusheat",
    "IH",
    "IndustryH",
    "SH"
)

$newModelTech = @(
    "standard_H",
    "BP",
    "standard_H",
    "BP",
    "standard_H",
    "standard_H",
    "BP",
    "standard_H",
    "standard_H",
    "BP",
    "HP",
    "standard_H",
    "HP",
    "HP",
    "HP",
    "standard_H",
    "standard_H"
)

$startRow = 7
for ($i = 0; $i -lt $newTech.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTech[$i]
    $ws.Cells.Item($row, 2).Value = $newModelTech[$i]
}
